# Generate Report for Archive
# The e9e86d34-... file moved back to "In Translation" (it is being retranslated),
# so it re-sorts ahead of 78170bc0-... (still "Ready for handoff") in every
# report sheet (Overview, zh-cn, de-de). This script swaps the row5/row6
# data + status and rebuilds the hyperlinks so display text follows the data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A5").Value = "e9e86d34-35fc-4917-ab78-af3bd2f61be7.md"
$ws.Range("B5").Value = "In Translation"
$ws.Range("C5").Value = "In Translation"
$ws.Range("D5").Value = "2016-30-18 00:30:39"

$ws.Range("A6").Value = "78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.md"
$ws.Range("B6").Value = "Ready for handoff"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "2016-28-18 00:28:46"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b236d27f6dd833f7f4a5e383568f97018d1d5d8a/e2e/70b6c8ce-f45c-48c3-b209-a11b72d043a5.md", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/103f0759b35a64be05c50900a17e297e40ff07ea/e2e/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md", "", "", "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/103f0759b35a64be05c50900a17e297e40ff07ea/e2e/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md", "", "", "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/1d66c2ad6e3d91fb6daeab0a1edea5db20de6b33/e2e/e9e86d34-35fc-4917-ab78-af3bd2f61be7.md", "", "", "e9e86d34-35fc-4917-ab78-af3bd2f61be7.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/dd1fb069f73170dceae4766754946ca7f12ca8c0/e2e/78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.md", "", "", "78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.md")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/d242c403c974b6952eab0bbe86aed0328aff2f8d/e2e/f43d7ba1-db92-40a5-b276-ec4111773384.md", "", "", "f43d7ba1-db92-40a5-b276-ec4111773384.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A5").Value = "e9e86d34-35fc-4917-ab78-af3bd2f61be7.md"
$ws.Range("C5").Value = "In Translation"
$ws.Range("D5").Value = "e9e86d34-35fc-4917-ab78-af3bd2f61be7.8877f16f869b27a8625bd997446535544be9b993.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-18 00:30:36"

$ws.Range("A6").Value = "78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.0c8fad67cd4dc2a987909bd9ea97f4545b629ff9.zh-cn.xlf"
$ws.Range("E6").Value = "2016-03-18 00:28:42"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b236d27f6dd833f7f4a5e383568f97018d1d5d8a/e2e/70b6c8ce-f45c-48c3-b209-a11b72d043a5.md", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/b236d27f6dd833f7f4a5e383568f97018d1d5d8a/e2e/70b6c8ce-f45c-48c3-b209-a11b72d043a5.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7571533ed125ae01df10791eb33567d14f922765/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.zh-cn.xlf", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a1f6a1ca491e5929d305b1a71a7c5ca70080f1f7/e2e/70b6c8ce-f45c-48c3-b209-a11b72d043a5.md", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/790e7ae7d55d41781720a6f8227bdcfadb8ebc55/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.zh-cn.xlf", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/103f0759b35a64be05c50900a17e297e40ff07ea/e2e/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md", "", "", "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/103f0759b35a64be05c50900a17e297e40ff07ea/e2e/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9dcfd52bc31ecb444f62cd03ca026efb7a965002/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.965bd93b8d7f0ff29d684f0a71d9bd2b1addfaa2.zh-cn.xlf", "", "", "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.965bd93b8d7f0ff29d684f0a71d9bd2b1addfaa2.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/103f0759b35a64be05c50900a17e297e40ff07ea/e2e/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md", "", "", "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/103f0759b35a64be05c50900a17e297e40ff07ea/e2e/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9dcfd52bc31ecb444f62cd03ca026efb7a965002/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.1669b18482b43a3dca05dabb6b7313f28bf459fc.zh-cn.xlf", "", "", "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.1669b18482b43a3dca05dabb6b7313f28bf459fc.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/1d66c2ad6e3d91fb6daeab0a1edea5db20de6b33/e2e/e9e86d34-35fc-4917-ab78-af3bd2f61be7.md", "", "", "e9e86d34-35fc-4917-ab78-af3bd2f61be7.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/1d66c2ad6e3d91fb6daeab0a1edea5db20de6b33/e2e/e9e86d34-35fc-4917-ab78-af3bd2f61be7.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d07c02d61c04d29121e7f1b32ea1bb83984ca680/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e9e86d34-35fc-4917-ab78-af3bd2f61be7.8877f16f869b27a8625bd997446535544be9b993.zh-cn.xlf", "", "", "e9e86d34-35fc-4917-ab78-af3bd2f61be7.8877f16f869b27a8625bd997446535544be9b993.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/dd1fb069f73170dceae4766754946ca7f12ca8c0/e2e/78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.md", "", "", "78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.md")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/dd1fb069f73170dceae4766754946ca7f12ca8c0/e2e/78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8f824b58beac68f8e35d24a1096863f5a89e867e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.0c8fad67cd4dc2a987909bd9ea97f4545b629ff9.zh-cn.xlf", "", "", "78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.0c8fad67cd4dc2a987909bd9ea97f4545b629ff9.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/d242c403c974b6952eab0bbe86aed0328aff2f8d/e2e/f43d7ba1-db92-40a5-b276-ec4111773384.md", "", "", "f43d7ba1-db92-40a5-b276-ec4111773384.md")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/d242c403c974b6952eab0bbe86aed0328aff2f8d/e2e/f43d7ba1-db92-40a5-b276-ec4111773384.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2c32c4b9f1224d5e6f729667faa826c1fb5340b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f43d7ba1-db92-40a5-b276-ec4111773384.f61309b95bd9958a707346f0cf3806e6dd42e3e7.zh-cn.xlf", "", "", "f43d7ba1-db92-40a5-b276-ec4111773384.f61309b95bd9958a707346f0cf3806e6dd42e3e7.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A5").Value = "e9e86d34-35fc-4917-ab78-af3bd2f61be7.md"
$ws.Range("C5").Value = "In Translation"
$ws.Range("D5").Value = "e9e86d34-35fc-4917-ab78-af3bd2f61be7.8877f16f869b27a8625bd997446535544be9b993.de-de.xlf"
$ws.Range("E5").Value = "2016-03-18 00:30:39"

$ws.Range("A6").Value = "78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.0c8fad67cd4dc2a987909bd9ea97f4545b629ff9.de-de.xlf"
$ws.Range("E6").Value = "2016-03-18 00:28:46"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b236d27f6dd833f7f4a5e383568f97018d1d5d8a/e2e/70b6c8ce-f45c-48c3-b209-a11b72d043a5.md", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/b236d27f6dd833f7f4a5e383568f97018d1d5d8a/e2e/70b6c8ce-f45c-48c3-b209-a11b72d043a5.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bad623942a2c732acb0ef0cc15d4a1bb697b0384/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.de-de.xlf", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/afcf07d6d226e8b675ee1d6ab7267ff1d09476af/e2e/70b6c8ce-f45c-48c3-b209-a11b72d043a5.md", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b37dbb7bea36555c450a1b8fb2d7b776c18fed74/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.de-de.xlf", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/103f0759b35a64be05c50900a17e297e40ff07ea/e2e/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md", "", "", "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/103f0759b35a64be05c50900a17e297e40ff07ea/e2e/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/22749cda8c148f07ea8b01186f9902b80873be0d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.965bd93b8d7f0ff29d684f0a71d9bd2b1addfaa2.de-de.xlf", "", "", "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.965bd93b8d7f0ff29d684f0a71d9bd2b1addfaa2.de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/103f0759b35a64be05c50900a17e297e40ff07ea/e2e/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md", "", "", "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/103f0759b35a64be05c50900a17e297e40ff07ea/e2e/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/22749cda8c148f07ea8b01186f9902b80873be0d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.1669b18482b43a3dca05dabb6b7313f28bf459fc.de-de.xlf", "", "", "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.1669b18482b43a3dca05dabb6b7313f28bf459fc.de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/1d66c2ad6e3d91fb6daeab0a1edea5db20de6b33/e2e/e9e86d34-35fc-4917-ab78-af3bd2f61be7.md", "", "", "e9e86d34-35fc-4917-ab78-af3bd2f61be7.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/1d66c2ad6e3d91fb6daeab0a1edea5db20de6b33/e2e/e9e86d34-35fc-4917-ab78-af3bd2f61be7.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f597a58d7f17547c9a650a23ac7f420b6f9dcc44/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e9e86d34-35fc-4917-ab78-af3bd2f61be7.8877f16f869b27a8625bd997446535544be9b993.de-de.xlf", "", "", "e9e86d34-35fc-4917-ab78-af3bd2f61be7.8877f16f869b27a8625bd997446535544be9b993.de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/dd1fb069f73170dceae4766754946ca7f12ca8c0/e2e/78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.md", "", "", "78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.md")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/dd1fb069f73170dceae4766754946ca7f12ca8c0/e2e/78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/54bea6992c8fdb2d1b7609ad859f3a51ee206247/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.0c8fad67cd4dc2a987909bd9ea97f4545b629ff9.de-de.xlf", "", "", "78170bc0-6a6e-40a1-ad2b-996d1f4fe20d.0c8fad67cd4dc2a987909bd9ea97f4545b629ff9.de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/d242c403c974b6952eab0bbe86aed0328aff2f8d/e2e/f43d7ba1-db92-40a5-b276-ec4111773384.md", "", "", "f43d7ba1-db92-40a5-b276-ec4111773384.md")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/d242c403c974b6952eab0bbe86aed0328aff2f8d/e2e/f43d7ba1-db92-40a5-b276-ec4111773384.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/73383fabe35d87a4557b4e92ef65f6b3f2e2dfb1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f43d7ba1-db92-40a5-b276-ec4111773384.f61309b95bd9958a707346f0cf3806e6dd42e3e7.de-de.xlf", "", "", "f43d7ba1-db92-40a5-b276-ec4111773384.f61309b95bd9958a707346f0cf3806e6dd42e3e7.de-de.xlf")

$wb.Save()
